$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5778633333333333
$ws.Range("H2").Value = 1.73359
$ws.Range("I2").Value = 0.01822222003587956
$ws.Range("J2").Value = 0.01829736853065394
$ws.Range("M2").Value = 4.277776333333333
$ws.Range("N2").Value = 12.833329
$ws.Range("O2").Value = 0.3536657835996513
$ws.Range("P2").Value = 0.3568846407551645
$ws.Range("Q2").Value = 2.471970091234444
$ws.Range("R2").Value = 22.24773082111
$ws.Range("S2").Value = 0.006444575727914612
$ws.Range("T2").Value = 0.006530049794827281
$ws.Range("G3").Value = 0.5778633333333333
$ws.Range("H3").Value = 1.73359
$ws.Range("I3").Value = 0.01822222003587956
$ws.Range("J3").Value = 0.01829736853065394
$ws.Range("M3").Value = 0.3272805
$ws.Range("N3").Value = 0.6545609999999999
$ws.Range("O3").Value = 0.02705796317293487
$ws.Range("P3").Value = 0.01820281918567982
$ws.Range("Q3").Value = 0.189123400665
$ws.Range("R3").Value = 1.13474040399
$ws.Range("S3").Value = 0.0004930561586599451
$ws.Range("T3").Value = 0.0003330636909372416
$ws.Range("G4").Value = 0.5778633333333333
$ws.Range("H4").Value = 1.73359
$ws.Range("I4").Value = 0.01822222003587956
$ws.Range("J4").Value = 0.01829736853065394
$ws.Range("M4").Value = 7.490476666666666
$ws.Range("N4").Value = 22.47143
$ws.Range("O4").Value = 0.6192762532274139
$ws.Range("P4").Value = 0.6249125400591558
$ws.Range("Q4").Value = 4.328471814855555
$ws.Range("R4").Value = 38.95624633369999
$ws.Range("S4").Value = 0.01128458814930501
$ws.Range("T4").Value = 0.01143425504488941
$ws.Range("I5").Value = 0.02935272232455987
$ws.Range("J5").Value = 0.02947377304702275
$ws.Range("M5").Value = 4.277776333333333
$ws.Range("N5").Value = 12.833329
$ws.Range("O5").Value = 0.3536657835996513
$ws.Range("P5").Value = 0.3568846407551645
$ws.Range("Q5").Value = 3.981899655461999
$ws.Range("R5").Value = 35.837096899158
$ws.Range("S5").Value = 0.01038105354169845
$ws.Range("T5").Value = 0.01051873690558596
$ws.Range("I6").Value = 0.02935272232455987
$ws.Range("J6").Value = 0.02947377304702275
$ws.Range("M6").Value = 0.3272805
$ws.Range("N6").Value = 0.6545609999999999
$ws.Range("O6").Value = 0.02705796317293487
$ws.Range("P6").Value = 0.01820281918567982
$ws.Range("Q6").Value = 0.3046438169369999
$ws.Range("R6").Value = 1.827862901622
$ws.Range("S6").Value = 0.0007942248796833242
$ws.Range("T6").Value = 0.0005365057614947185
$ws.Range("I7").Value = 0.02935272232455987
$ws.Range("J7").Value = 0.02947377304702275
$ws.Range("M7").Value = 7.490476666666666
$ws.Range("N7").Value = 22.47143
$ws.Range("O7").Value = 0.6192762532274139
$ws.Range("P7").Value = 0.6249125400591558
$ws.Range("Q7").Value = 6.972390357539999
$ws.Range("R7").Value = 62.75151321785999
$ws.Range("S7").Value = 0.0181774439031781
$ws.Range("T7").Value = 0.01841853037994207
$ws.Range("G8").Value = 16.208374
$ws.Range("H8").Value = 48.625122
$ws.Range("I8").Value = 0.5111114348580046
$ws.Range("J8").Value = 0.5132192600799546
$ws.Range("M8").Value = 4.277776333333333
$ws.Range("N8").Value = 12.833329
$ws.Range("O8").Value = 0.3536657835996513
$ws.Range("P8").Value = 0.3568846407551645
$ws.Range("Q8").Value = 69.33579869901534
$ws.Range("R8").Value = 624.022188291138
$ws.Range("S8").Value = 0.1807626261157984
$ws.Range("T8").Value = 0.1831600712622659
$ws.Range("G9").Value = 16.208374
$ws.Range("H9").Value = 48.625122
$ws.Range("I9").Value = 0.5111114348580046
$ws.Range("J9").Value = 0.5132192600799546
$ws.Range("M9").Value = 0.3272805
$ws.Range("N9").Value = 0.6545609999999999
$ws.Range("O9").Value = 0.02705796317293487
$ws.Range("P9").Value = 0.01820281918567982
$ws.Range("Q9").Value = 5.304684746907
$ws.Range("R9").Value = 31.828108481442
$ws.Range("S9").Value = 0.01382963438165379
$ws.Range("T9").Value = 0.009342037393843797
$ws.Range("G10").Value = 16.208374
$ws.Range("H10").Value = 48.625122
$ws.Range("I10").Value = 0.5111114348580046
$ws.Range("J10").Value = 0.5132192600799546
$ws.Range("M10").Value = 7.490476666666666
$ws.Range("N10").Value = 22.47143
$ws.Range("O10").Value = 0.6192762532274139
$ws.Range("P10").Value = 0.6249125400591558
$ws.Range("Q10").Value = 121.4084472516067
$ws.Range("R10").Value = 1092.67602526446
$ws.Range("S10").Value = 0.3165191743605525
$ws.Range("T10").Value = 0.3207171514238449
$ws.Range("G11").Value = 0.39073
$ws.Range("H11").Value = 0.78146
$ws.Range("I11").Value = 0.01232119711342224
$ws.Range("J11").Value = 0.008248006513630573
$ws.Range("M11").Value = 4.277776333333333
$ws.Range("N11").Value = 12.833329
$ws.Range("O11").Value = 0.3536657835996513
$ws.Range("P11").Value = 0.3568846407551645
$ws.Range("Q11").Value = 1.671455546723333
$ws.Range("R11").Value = 10.02873328034
$ws.Range("S11").Value = 0.00435758583200424
$ws.Range("T11").Value = 0.002943586841563303
$ws.Range("G12").Value = 0.39073
$ws.Range("H12").Value = 0.78146
$ws.Range("I12").Value = 0.01232119711342224
$ws.Range("J12").Value = 0.008248006513630573
$ws.Range("M12").Value = 0.3272805
$ws.Range("N12").Value = 0.6545609999999999
$ws.Range("O12").Value = 0.02705796317293487
$ws.Range("P12").Value = 0.01820281918567982
$ws.Range("Q12").Value = 0.127878309765
$ws.Range("R12").Value = 0.51151323906
$ws.Range("S12").Value = 0.0003333864977414505
$ws.Range("T12").Value = 0.0001501369712099267
$ws.Range("G13").Value = 0.39073
$ws.Range("H13").Value = 0.78146
$ws.Range("I13").Value = 0.01232119711342224
$ws.Range("J13").Value = 0.008248006513630573
$ws.Range("M13").Value = 7.490476666666666
$ws.Range("N13").Value = 22.47143
$ws.Range("O13").Value = 0.6192762532274139
$ws.Range("P13").Value = 0.6249125400591558
$ws.Range("Q13").Value = 2.926753947966667
$ws.Range("R13").Value = 17.5605236878
$ws.Range("S13").Value = 0.007630224783676554
$ws.Range("T13").Value = 0.005154282700857343
$ws.Range("G14").Value = 13.60421466666667
$ws.Range("H14").Value = 40.81264400000001
$ws.Range("I14").Value = 0.4289924256681337
$ws.Range("J14").Value = 0.4307615918287382
$ws.Range("M14").Value = 4.277776333333333
$ws.Range("N14").Value = 12.833329
$ws.Range("O14").Value = 0.3536657835996513
$ws.Range("P14").Value = 0.3568846407551645
$ws.Range("Q14").Value = 58.19578753465289
$ws.Range("R14").Value = 523.7620878118761
$ws.Range("S14").Value = 0.1517199423822357
$ws.Range("T14").Value = 0.153732195950922
$ws.Range("G15").Value = 13.60421466666667
$ws.Range("H15").Value = 40.81264400000001
$ws.Range("I15").Value = 0.4289924256681337
$ws.Range("J15").Value = 0.4307615918287382
$ws.Range("M15").Value = 0.3272805
$ws.Range("N15").Value = 0.6545609999999999
$ws.Range("O15").Value = 0.02705796317293487
$ws.Range("P15").Value = 0.01820281918567982
$ws.Range("Q15").Value = 4.452394178214001
$ws.Range("R15").Value = 26.714365069284
$ws.Range("S15").Value = 0.01160766125519636
$ws.Range("T15").Value = 0.007841075368194134
$ws.Range("G16").Value = 13.60421466666667
$ws.Range("H16").Value = 40.81264400000001
$ws.Range("I16").Value = 0.4289924256681337
$ws.Range("J16").Value = 0.4307615918287382
$ws.Range("M16").Value = 7.490476666666666
$ws.Range("N16").Value = 22.47143
$ws.Range("O16").Value = 0.6192762532274139
$ws.Range("P16").Value = 0.6249125400591558
$ws.Range("Q16").Value = 101.9020525289911
$ws.Range("R16").Value = 917.1184727609201
$ws.Range("S16").Value = 0.2656648220307017
$ws.Range("T16").Value = 0.269188320509622
